# SS : Add crypto address in ID
# A new "Crypto Address" column is inserted right before the existing
# "Succession" column (column N), pushing Succession / Newsletter /
# Password / Tries one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at N (14) - mirrors right-click > Insert on the
# column header; existing formatting/data in N..Q shifts to O..R.
$ws.Columns.Item(14).Insert()

# Fill in the new "Crypto Address" column header + the sample row value.
$ws.Cells.Item(1, 14).Value = "Crypto Address"
$ws.Cells.Item(2, 14).Value = "CCLcWAJX6fubUqGyZWz8dyUGEddRj8h4XZZCNSD"

# Restore the frozen panes (2 columns / 1 row) that existed before the
# column insert, keeping the same split point.
$ws.Range("C2").Select() | Out-Null
$win = $excel.ActiveWindow
$win.FreezePanes = $true

# Leave the selection where the user ended up after adding the column.
$ws.Range("N3").Select() | Out-Null
